$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark from its old location
#    (it currently sits just before the "Install ERP package using the
#    latest release URL:" run).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Trim the trailing "(Search customer orders on CRM to get latest
#    info)" remark off the "... field values:" sentence, then append the
#    new license-count guidance as its own run right after it.
# ------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute(": (Search customer orders on CRM to get latest info)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$newText = " based on SSO or Sales Order in the CRM org.  Do NOT include RSF licenses in this count."

$insertPoint = $d.Range($r.End, $r.End)
$insertPoint.InsertAfter($newText)

$newRunEnd = $r.End + $newText.Length
$newRun = $d.Range($r.End, $newRunEnd)
$newRun.Font.Color = 12611584

# ------------------------------------------------------------------
# 3. Re-add the "_GoBack" bookmark collapsed at the end of this same
#    paragraph (right after the text we just inserted, before the
#    paragraph mark). A temporary placeholder character is inserted and
#    then removed around the target point - adding a collapsed bookmark
#    exactly at a paragraph's last position is mishandled otherwise.
# ------------------------------------------------------------------
$placeholderPoint = $d.Range($newRunEnd, $newRunEnd)
$placeholderPoint.InsertAfter("X")

$bmPos = $d.Range($newRunEnd, $newRunEnd)
$d.Bookmarks.Add("_GoBack", $bmPos)

$placeholder = $d.Range($newRunEnd, $newRunEnd + 1)
$placeholder.Delete()
